# Preenche a linha 8 com o resultado da leitura (RESPOSTA) e o valor lido.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "RESPOSTA "
$ws.Range("B8").Value = "000000400000000000a34a47002f8b56a0afc105401e441d0000000000000089000000400000000080a5283c403c712300d8161980534d210000000000000034"

$ws.Rows.Item(8).RowHeight = 12.85

$null = $ws.Range("C11").Select()
